$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.455362044514542
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 25.99591228164478

$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 27.82738278199502

$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.7527432677738641
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 6.189590430959694

$ws.Range("B5").Value = 3.286832544864788
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("G5").Value = 6.189590430959694
